$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C26").Value = "Land"
$ws.Range("E26").Value = "Normal"

$ws.Range("E26").Select()
